$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.113.40'
$ws.Cells.Item(2, 5).Value = '  -0.09%  '

$ws.Cells.Item(3, 4).Value = '1.655.97'
$ws.Cells.Item(3, 5).Value = '  -0.10%  '

$ws.Cells.Item(4, 5).Value = '  -0.20%  '

$ws.Cells.Item(5, 4).Value = '217.81'
$ws.Cells.Item(5, 5).Value = '  +0.59%  '

$ws.Cells.Item(6, 4).Value = '0.5262'
$ws.Cells.Item(6, 5).Value = '  +1.85%  '

$ws.Cells.Item(7, 5).Value = '  -0.17%  '

$ws.Cells.Item(8, 4).Value = '0.2613'
$ws.Cells.Item(8, 5).Value = '  -0.71%  '

$ws.Cells.Item(9, 4).Value = '0.06362'
$ws.Cells.Item(9, 5).Value = '  +1.48%  '

$ws.Cells.Item(10, 4).Value = '20.47'
$ws.Cells.Item(10, 5).Value = '  -1.37%  '

$ws.Cells.Item(11, 4).Value = '0.07807'
$ws.Cells.Item(11, 5).Value = '  +1.22%  '

$ws.Cells.Item(12, 4).Value = '4.515'
$ws.Cells.Item(12, 5).Value = '  +2.02%  '

$ws.Cells.Item(13, 4).Value = '1.658.39'
$ws.Cells.Item(13, 5).Value = '  -0.01%  '

$ws.Cells.Item(14, 4).Value = '0.5497'
$ws.Cells.Item(14, 5).Value = '  +1.52%  '

$ws.Cells.Item(15, 4).Value = '0.0₅8237'
$ws.Cells.Item(15, 5).Value = '  +1.50%  '

$ws.Cells.Item(16, 4).Value = '65.49'
$ws.Cells.Item(16, 5).Value = '  +1.19%  '

$ws.Cells.Item(17, 4).Value = '26.116.02'
$ws.Cells.Item(17, 5).Value = '  -0.17%  '

$ws.Cells.Item(18, 5).Value = '  -0.22%  '

$ws.Cells.Item(19, 4).Value = '4.602'
$ws.Cells.Item(19, 5).Value = '  -0.26%  '

$ws.Cells.Item(20, 4).Value = '191.39'
$ws.Cells.Item(20, 5).Value = '  -0.16%  '

$ws.Cells.Item(21, 5).Value = '  -0.04%  '

$ws.Cells.Item(22, 4).Value = '6.047'
$ws.Cells.Item(22, 5).Value = '  +0.53%  '

$ws.Cells.Item(23, 4).Value = '1.003'
$ws.Cells.Item(23, 5).Value = '  -0.22%  '

$ws.Cells.Item(24, 4).Value = '141.91'
$ws.Cells.Item(24, 5).Value = '  +1.52%  '

$ws.Cells.Item(25, 4).Value = '0.1238'
$ws.Cells.Item(25, 5).Value = '  +1.03%  '

$ws.Cells.Item(26, 4).Value = '7.258'
$ws.Cells.Item(26, 5).Value = '  +1.14%  '

$ws.Cells.Item(27, 4).Value = '16.12'
$ws.Cells.Item(27, 5).Value = '  +0.19%  '

$ws.Cells.Item(28, 4).Value = '1.427'
$ws.Cells.Item(28, 5).Value = '  +1.54%  '

$ws.Cells.Item(29, 4).Value = '0.05904'
$ws.Cells.Item(29, 5).Value = '  -1.28%  '

$ws.Cells.Item(30, 5).Value = '  +0.31%  '

$ws.Cells.Item(31, 4).Value = '3.516'
$ws.Cells.Item(31, 5).Value = '  -0.88%  '

$ws.Cells.Item(32, 4).Value = '3.264'
$ws.Cells.Item(32, 5).Value = '  +0.31%  '

$ws.Cells.Item(33, 5).Value = '  -0.38%  '

$ws.Cells.Item(34, 4).Value = '0.9534'
$ws.Cells.Item(34, 5).Value = '  -1.15%  '

$ws.Cells.Item(35, 5).Value = '  +0.54%  '

$ws.Cells.Item(36, 4).Value = '2.412'
$ws.Cells.Item(36, 5).Value = '  -0.49%  '

$ws.Cells.Item(37, 4).Value = '0.5708'
$ws.Cells.Item(37, 5).Value = '  +0.36%  '

$ws.Cells.Item(38, 4).Value = '0.01622'
$ws.Cells.Item(38, 5).Value = '  +2.00%  '

$ws.Cells.Item(39, 4).Value = '5.821'
$ws.Cells.Item(39, 5).Value = '  -2.37%  '

$ws.Cells.Item(40, 5).Value = '  -0.74%  '

$ws.Cells.Item(41, 5).Value = '  -0.08%  '

$ws.Cells.Item(42, 4).Value = '1.031.48'
$ws.Cells.Item(42, 5).Value = '  +2.50%  '

$ws.Cells.Item(43, 4).Value = '102.83'
$ws.Cells.Item(43, 5).Value = '  +2.50%  '

$ws.Cells.Item(44, 4).Value = '1.798.43'
$ws.Cells.Item(44, 5).Value = '  -0.01%  '

$ws.Cells.Item(45, 4).Value = '57.27'
$ws.Cells.Item(45, 5).Value = '  +1.04%  '

$ws.Cells.Item(46, 4).Value = '1.003'
$ws.Cells.Item(46, 5).Value = '  -0.32%  '

$ws.Cells.Item(47, 4).Value = '0.4298'
$ws.Cells.Item(47, 5).Value = '  +2.81%  '

$ws.Cells.Item(48, 4).Value = '1.479'
$ws.Cells.Item(48, 5).Value = '  +2.29%  '

$ws.Cells.Item(49, 4).Value = '7.858'
$ws.Cells.Item(49, 5).Value = '  -1.70%  '

$ws.Cells.Item(50, 4).Value = '0.05151'
$ws.Cells.Item(50, 5).Value = '  -0.33%  '

$ws.Cells.Item(51, 4).Value = '0.09716'
$ws.Cells.Item(51, 5).Value = '  -0.04%  '
